$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above down into the new row, then fill in values
$ws.Range("A8:F8").Copy() | Out-Null
$ws.Range("A9:F9").PasteSpecial(-4122) | Out-Null

$ws.Range("A9").Value = "EXL_CorporateLensHomePage_MarketPlace"
$ws.Range("B9").Value = "Add New Marketplace"
$ws.Range("C9").Value = "N"
$ws.Range("D9").Value = "Y"
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = "Sprint1"

# Extend the existing data validations to include the new row
$ws.Range("C2:D8").Validation.Delete()
$ws.Range("F2:F8").Validation.Delete()

$ws.Range("C2:D9").Validation.Add(3, 1, 1, """Y,N""")
$ws.Range("C2:D9").Validation.IgnoreBlank = $true
$ws.Range("C2:D9").Validation.InCellDropdown = $true
$ws.Range("C2:D9").Validation.ShowInput = $true
$ws.Range("C2:D9").Validation.ShowError = $true

$ws.Range("F2:F9").Validation.Add(3, 1, 1, """Sprint1,Sprint2,Sprint3,Sprint4,Sprint5,Sprint6,Sprint7,Sprint8,Sprint9,Sprint10""")
$ws.Range("F2:F9").Validation.IgnoreBlank = $true
$ws.Range("F2:F9").Validation.InCellDropdown = $true
$ws.Range("F2:F9").Validation.ShowInput = $true
$ws.Range("F2:F9").Validation.ShowError = $true
